$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume text columns keep their literal text representation
# (values like "1.001" or "0.5200" would otherwise be re-interpreted as numbers)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "25.982.13"
$ws.Cells.Item(2, 5).Value = "  +0.70%  "
$ws.Cells.Item(3, 4).Value = "1.748.05"
$ws.Cells.Item(4, 5).Value = "  +0.11%  "
$ws.Cells.Item(5, 4).Value = "234.06"
$ws.Cells.Item(5, 5).Value = "  -1.26%  "
$ws.Cells.Item(6, 4).Value = "1.001"
$ws.Cells.Item(6, 5).Value = "  +0.13%  "
$ws.Cells.Item(7, 4).Value = "0.5194"
$ws.Cells.Item(7, 5).Value = "  +2.03%  "
$ws.Cells.Item(8, 4).Value = "0.2822"
$ws.Cells.Item(8, 5).Value = "  +4.48%  "
$ws.Cells.Item(9, 4).Value = "39.55"
$ws.Cells.Item(9, 5).Value = "  -3.72%  "
$ws.Cells.Item(10, 4).Value = "0.06136"
$ws.Cells.Item(10, 5).Value = "  -1.15%  "
$ws.Cells.Item(11, 4).Value = "1.751.22"
$ws.Cells.Item(11, 5).Value = "  -0.11%  "
$ws.Cells.Item(12, 4).Value = "0.07018"
$ws.Cells.Item(12, 5).Value = "  +1.45%  "
$ws.Cells.Item(13, 4).Value = "15.42"
$ws.Cells.Item(13, 5).Value = "  -0.97%  "
$ws.Cells.Item(14, 4).Value = "0.6439"
$ws.Cells.Item(14, 5).Value = "  +2.75%  "
$ws.Cells.Item(15, 4).Value = "4.531"
$ws.Cells.Item(15, 5).Value = "  +0.77%  "
$ws.Cells.Item(16, 4).Value = "77.50"
$ws.Cells.Item(16, 5).Value = "  -1.50%  "
$ws.Cells.Item(17, 5).Value = "  +0.12%  "
$ws.Cells.Item(18, 4).Value = "1.001"
$ws.Cells.Item(18, 5).Value = "  +0.08%  "
$ws.Cells.Item(19, 4).Value = "25.983.46"
$ws.Cells.Item(19, 5).Value = "  +0.61%  "
$ws.Cells.Item(20, 5).Value = "  -1.77%  "
$ws.Cells.Item(21, 5).Value = "  -1.56%  "
$ws.Cells.Item(22, 4).Value = "1.981.20"
$ws.Cells.Item(22, 5).Value = "  +0.38%  "
$ws.Cells.Item(23, 4).Value = "4.155"
$ws.Cells.Item(23, 5).Value = "  +1.96%  "
$ws.Cells.Item(24, 4).Value = "8.653"
$ws.Cells.Item(24, 5).Value = "  +4.65%  "
$ws.Cells.Item(25, 4).Value = "5.153"
$ws.Cells.Item(25, 5).Value = "  -0.59%  "
$ws.Cells.Item(26, 4).Value = "139.16"
$ws.Cells.Item(26, 5).Value = "  +1.71%  "
$ws.Cells.Item(27, 4).Value = "1.508"
$ws.Cells.Item(27, 5).Value = "  +2.96%  "
$ws.Cells.Item(28, 4).Value = "1.828"
$ws.Cells.Item(28, 5).Value = "  +1.53%  "
$ws.Cells.Item(29, 5).Value = "  -1.51%  "
$ws.Cells.Item(30, 4).Value = "102.83"
$ws.Cells.Item(30, 5).Value = "  +0.06%  "
$ws.Cells.Item(31, 4).Value = "0.08293"
$ws.Cells.Item(31, 5).Value = "  +0.32%  "
$ws.Cells.Item(32, 4).Value = "3.670"
$ws.Cells.Item(32, 5).Value = "  -1.82%  "
$ws.Cells.Item(33, 4).Value = "3.438"
$ws.Cells.Item(33, 5).Value = "  -0.09%  "
$ws.Cells.Item(34, 4).Value = "0.04482"
$ws.Cells.Item(34, 5).Value = "  +1.56%  "
$ws.Cells.Item(35, 4).Value = "2.617"
$ws.Cells.Item(35, 5).Value = "  -0.91%  "
$ws.Cells.Item(36, 4).Value = "0.9896"
$ws.Cells.Item(36, 5).Value = "  -1.73%  "
$ws.Cells.Item(37, 4).Value = "0.6152"
$ws.Cells.Item(37, 5).Value = "  +1.48%  "
$ws.Cells.Item(38, 4).Value = "2.677"
$ws.Cells.Item(38, 5).Value = "  -0.63%  "
$ws.Cells.Item(39, 4).Value = "0.01586"
$ws.Cells.Item(39, 5).Value = "  +1.66%  "
$ws.Cells.Item(40, 4).Value = "1.932"
$ws.Cells.Item(40, 5).Value = "  -2.02%  "
$ws.Cells.Item(41, 4).Value = "1.001"
$ws.Cells.Item(41, 5).Value = "  +0.07%  "
$ws.Cells.Item(42, 2).Value = "Quant"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(42, 4).Value = "100.66"
$ws.Cells.Item(42, 5).Value = "  -1.53%  "
$ws.Cells.Item(43, 2).Value = "TheSandbox"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(43, 4).Value = "0.3869"
$ws.Cells.Item(43, 5).Value = "  +0.26%  "
$ws.Cells.Item(44, 2).Value = "TrustWalletToken"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(44, 4).Value = "0.7395"
$ws.Cells.Item(44, 5).Value = "  -1.47%  "
$ws.Cells.Item(45, 4).Value = "5.066"
$ws.Cells.Item(45, 5).Value = "  +3.86%  "
$ws.Cells.Item(46, 2).Value = "Cronos"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(46, 4).Value = "0.05469"
$ws.Cells.Item(46, 5).Value = "  -0.71%  "
$ws.Cells.Item(47, 2).Value = "Aptos"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(47, 4).Value = "6.309"
$ws.Cells.Item(47, 5).Value = "  +5.67%  "
$ws.Cells.Item(48, 2).Value = "Algorand"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(48, 4).Value = "0.1126"
$ws.Cells.Item(48, 5).Value = "  +2.76%  "
$ws.Cells.Item(49, 2).Value = "Aave"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(49, 4).Value = "53.00"
$ws.Cells.Item(49, 5).Value = "  +0.08%  "
$ws.Cells.Item(50, 2).Value = "Elrond"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(50, 4).Value = "30.02"
$ws.Cells.Item(50, 5).Value = "  -0.97%  "
$ws.Cells.Item(51, 2).Value = "EnergySwap"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(51, 4).Value = "7.659"
$ws.Cells.Item(51, 5).Value = "  +2.36%  "
